$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 186: Ajo Chino "1a (cosecha)"
$ws.Cells.Item(186, 1).Value2 = 9
$ws.Cells.Item(186, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186, 3).Value2 = "Metropolitana"
$ws.Cells.Item(186, 4).Value2 = 44595
$ws.Cells.Item(186, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(186, 5).Value2 = 13
$ws.Cells.Item(186, 6).Value2 = 100112003
$ws.Cells.Item(186, 7).Value2 = "Ajo"
$ws.Cells.Item(186, 8).Value2 = "Chino"
$ws.Cells.Item(186, 9).Value2 = "1a (cosecha)"
$ws.Cells.Item(186, 10).Value2 = 1600
$ws.Cells.Item(186, 11).Value2 = 3000
$ws.Cells.Item(186, 12).Value2 = 3000
$ws.Cells.Item(186, 13).Value2 = 3000
$ws.Cells.Item(186, 14).Value2 = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(186, 15).Value2 = "Provincia de Talagante"
$ws.Cells.Item(186, 16).Value2 = 150
$ws.Cells.Item(186, 17).Value2 = 20
$ws.Cells.Item(186, 18).Value2 = "Hortaliza"

# New row 187: Ajo Chino "2a (cosecha)"
$ws.Cells.Item(187, 1).Value2 = 9
$ws.Cells.Item(187, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(187, 3).Value2 = "Metropolitana"
$ws.Cells.Item(187, 4).Value2 = 44595
$ws.Cells.Item(187, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(187, 5).Value2 = 13
$ws.Cells.Item(187, 6).Value2 = 100112003
$ws.Cells.Item(187, 7).Value2 = "Ajo"
$ws.Cells.Item(187, 8).Value2 = "Chino"
$ws.Cells.Item(187, 9).Value2 = "2a (cosecha)"
$ws.Cells.Item(187, 10).Value2 = 800
$ws.Cells.Item(187, 11).Value2 = 2000
$ws.Cells.Item(187, 12).Value2 = 2000
$ws.Cells.Item(187, 13).Value2 = 2000
$ws.Cells.Item(187, 14).Value2 = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(187, 15).Value2 = "Provincia de Talagante"
$ws.Cells.Item(187, 16).Value2 = 100
$ws.Cells.Item(187, 17).Value2 = 20
$ws.Cells.Item(187, 18).Value2 = "Hortaliza"
